$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value()
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 5.95 = 24250.0 pesos"), "✅ 1000 Bs = 5.56 = 22537.84 pesos"
$text = $text -replace [regex]::Escape("✅ 24250.0 pesos = 5.92 = 959.34 Bs"), "✅ 22537.84 pesos = 5.53 = 926.74 Bs"
$cell.Value = $text

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 179.7
$ws2.Range("O10").Value = 4050.05
$ws2.Range("N12").Value = 4073.5
$ws2.Range("O12").Value = 167.5
